$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New building-block rows for d0022-d0024 (rows 72-74)
$ws.Range("A72").Value = 'd0022'
$ws.Range("C72").Value = '$x^{3}+x^{2}+C$;'
$ws.Range("B72").Value = '부정적분을 이용해서 적분상수가 표함된 함수 $f(x)$ 의 꼴을 구합니다.'
$ws.Range("A73").Value = 'd0023'
$ws.Range("B73").Value = '주어진 함숫값을 이용해서 적분상수값을 정해주고 함수 $f(x)$를 결정해 줍니다.'
$ws.Range("A74").Value = 'd0024'
$ws.Range("B74").Value = '문제에서 요구하는 함숫값을 구합니다.'
$ws.Range("C74").Value = '$f(1)$;'

# New common-part problems 17 & 18 (rows 42-43)
$ws.Range("A42").Value = 'c0032'
$ws.Range("B42").Value = '$\displaystyle\sum$의 성질을 이용해서 $\displaystyle\sum$가 포함된 식을 정리합니다.'
$ws.Range("A43").Value = 'c0033'
$ws.Range("B43").Value = '$\displaystyle\sum$의 성질을 이용해서 $\displaystyle\sum$가 포함된 두 식을 연립합니다.'

# Update C1 label to include the new $a_{8}$ term
$ws.Range("C1").Value = '$a_{10}$; $k$; $a_{8}$;'

# Move the view/selection to C1 (top-left scroll anchor + active cell)
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C1").Select()
